$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts old B->C, old C->D) to make
# room for the new "StatQuery" column.
$ws.Columns.Item(2).Insert()

# New header + stat-bar query text in the freshly inserted column B.
$ws.Range("B1").Value = "StatQuery"
$ws.Range("B2").Value = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['German Shorthaired Pointer']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

# Match column A's width on the new column B so both share the wide
# "query text" column formatting.
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(1).ColumnWidth

# Note: column Insert() already copies A2's wrap-text cell format onto the
# new B2 cell, so no extra style assignment is needed here (re-assigning
# .Style from another cell's .Style would reset it to the default "Normal"
# named style and lose the wrap-text direct formatting).

# Update the active selection to match the refreshed sheet (whole column B).
[void]$ws.Range("B1:B1048576").Select()

Write-Output "edit complete"
